# CRM-2487: Add "Parts Code" column to the Spare Requested Parts export
# template, between "Age of Requested" and "Parts Required".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at F (shifts Parts Required .. Brand one column right).
$ws.Columns("F:F").Insert()

# Populate the new header (row 1) and its merge-placeholder (row 2).
$ws.Range("F1").Value = "Parts Code"
$ws.Range("F2").Value = "{spare:part_number}"
